$d = $word.ActiveDocument

# Update the date heading paragraph.
$d.Content.Find.Execute("2024-05-13 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-05-14 Tuesday", 2)

# Update the answer table. It has 5 columns; only every 4th row (1, 5, 9, 13, 17)
# holds the worked-answer text, the rows between them are blank spacer rows.
$t = $d.Tables.Item(1)

$rowUpdates = @(
    @{ Row = 1;  Values = @("12÷2=6, 0", "46÷5=9, 1", "58÷2=29, 0", "69÷9=7, 6", "71÷2=35, 1") },
    @{ Row = 5;  Values = @("64÷4=16, 0", "40÷5=8, 0", "57÷7=8, 1", "11÷5=2, 1", "87÷4=21, 3") },
    @{ Row = 9;  Values = @("51÷5=10, 1", "77÷5=15, 2", "48÷2=24, 0", "99÷7=14, 1", "41÷8=5, 1") },
    @{ Row = 13; Values = @("36÷4=9, 0", "84÷9=9, 3", "82÷5=16, 2", "12÷7=1, 5", "53÷7=7, 4") },
    @{ Row = 17; Values = @("23÷8=2, 7", "22÷2=11, 0", "43÷3=14, 1", "58÷7=8, 2", "32÷9=3, 5") }
)

foreach ($update in $rowUpdates) {
    $rowIndex = $update.Row
    $values = $update.Values
    for ($col = 1; $col -le $values.Length; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
